# Swap the B, E, F, G column values between each pair of rows listed below.
# These pairs represent duplicate line items (same item code/description/rate)
# whose sale-return and sale rows were mixed up; the fix swaps the
# code/price/qty/value figures back between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(142,143),
    @(154,155),
    @(176,177),
    @(256,257),
    @(305,306),
    @(308,310),
    @(338,339),
    @(343,344),
    @(364,365),
    @(367,368),
    @(371,372),
    @(392,393),
    @(449,450),
    @(582,583),
    @(591,592),
    @(701,702),
    @(707,708)
)

$cols = @("B", "E", "F", "G")

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    foreach ($col in $cols) {
        $cellA = $ws.Range("$col$rowA")
        $cellB = $ws.Range("$col$rowB")

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}
